$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Activation date: 2012 -> 2024
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Ativação: 01/01/2012", $false, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2024", 2)

# ---------------------------------------------------------------------
# 2. Curso (semestre ideal): drop the "EQD (7), " part
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Curso (semestre ideal): EQD (7), EQN (8)", $false, $false, $false, $false, $false, $true, 1, $false, "Curso (semestre ideal): EQN (8)", 2)

# ---------------------------------------------------------------------
# 3. Insert an italic English translation paragraph right after the
#    "Objetivos" (Capacitar os alunos...) paragraph.
# ---------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("Capacitar os alunos a calcular os parâmetros de projeto de reatores ideais, a distinguir entre um reator ideal e um real, e a compreender a influência da temperatura e pressão no projeto de reatores químicos.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$objPara = $findRange.Paragraphs(1)
$objPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $objPara.Next()
$newRange = $newPara.Range
$newRange.Text = "Enable students to calculate the project parameters of ideal reactors, distinguish between an ideal and a real reactor, and understand the influence of temperature and pressure on the design of chemical reactors."
$italicRange = $objPara.Next().Range
$italicRange.MoveEnd(1, -1)
$italicRange.Font.Italic = $true

# ---------------------------------------------------------------------
# 4. Docentes: add a second professor on a new line (w:br) inside the
#    same list-bullet paragraph, as a brand-new run.
# ---------------------------------------------------------------------
$findRange2 = $d.Content
$findRange2.Find.Execute("5963230 - Leandro Gonçalves de Aguiar", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$docPara = $findRange2.Paragraphs(1)
$docPara.Range.InsertParagraphAfter() | Out-Null
$newDocPara = $docPara.Next()
$newDocPara.Range.Text = "6310316 - Liana Alvares Rodrigues"
# Turn the boundary between the two paragraphs into a line break so both
# pieces of text end up inside the *same* paragraph, as two runs.
$markRange = $d.Range($docPara.Range.End - 1, $docPara.Range.End)
$markRange.InsertBefore([char]11)
$docPara2 = $findRange2.Paragraphs(1)
$markRange2 = $d.Range($docPara2.Range.End - 1, $docPara2.Range.End)
$markRange2.Delete()

# ---------------------------------------------------------------------
# 5. Insert an italic English translation paragraph right after the
#    "Programa resumido" content paragraph.
# ---------------------------------------------------------------------
$findRange3 = $d.Content
$findRange3.Find.Execute("1. Introdução a Reatores. 2. Modelos Ideais de Reatores Químicos Isotérmicos  Reações Simples. 3. Reações Múltiplas em Reatores Ideais. 4. Efeitos Térmicos em Reatores Ideais. 5. Reatores Catalíticos Heterogêneos. 6. Reatores Não-Ideais", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$resumoPara = $findRange3.Paragraphs(1)
$resumoPara.Range.InsertParagraphAfter() | Out-Null
$newResumoPara = $resumoPara.Next()
$newResumoPara.Range.Text = "1. Introduction to Reactors, 2. Ideal Models of Isothermal Chemical Reactors - Simple Reactions, 3. Multiple Reactions in Ideal Reactors, 4. Thermal Effects in Ideal Reactors, 5. Heterogeneous Catalytic Reactors, 6. Non-Ideal Reactors."
$italicRange2 = $resumoPara.Next().Range
$italicRange2.MoveEnd(1, -1)
$italicRange2.Font.Italic = $true

# ---------------------------------------------------------------------
# 6. Programa: collapse the whole bulleted/line-broken list into a
#    single run of text (no more w:br elements).
# ---------------------------------------------------------------------
$findRange4 = $d.Content
$findRange4.Find.Execute("1. Introdução a Reatores: Conceitos básicos", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$programaPara = $findRange4.Paragraphs(1)
$programaPara.Range.Text = "1. Introdução a Reatores: Conceitos básicos2. Modelos Ideais de Reatores Químicos Isotérmicos  Reações Simples: 2.1) Equações fundamentais de projeto de reatores2.2) Reator tanque descontínuo (BSTR)2.3) Reator tanque de mistura contínuo (CSTR)2.4) Reator tubular de fluxo pistonado (PFR)2.5) Comparação de desempenho de reatores CSTR e PFR2.6) Reatores CSTR em cascata2.7) Associação mista de reatores em série: CSTR e PFR2.8) Reatores com reciclo2.9) Reações auto-catalíticas2.10) Reatores semi-contínuos3. Reações Múltiplas em Reatores Ideais3.1) Noções gerais: otimização, rendimento e seletividade3.2) Reações paralelas e reações em série3.3) Sistemas com reações série-paralelo: reações de múltipla substituição e reações poliméricas3.4) Problemas simples de otimização4. Efeitos Térmicos em Reatores Ideais4.1) Equação do balanço de energia4.2) Balanço de energia aplicado ao BSTR4.3) Balanço de energia aplicado ao CSTR4.4) Balanço de energia aplicado ao PFR5. Reatores Catalíticos Heterogêneos5.1) Introdução5.2) Efeito dos processos físicos sobre a taxa de reação5.2.1  Fenômenos interfases5.2.2  Fenômenos intrapartícula5.2.3  Difusão e reação em catalisadores porosos5.3) Cálculo de reatores de leito fixo5.4) Reatores trifásicos6. Reatores Não-Ideais6.1) A distribuição dos tempos de residência6.2) Modelos dos tanques contínuos em série6.3) Modelo da dispersão axial"

# ---------------------------------------------------------------------
# 7. Insert an italic English translation paragraph right after the
#    (now single-run) "Programa" paragraph.
# ---------------------------------------------------------------------
$findRange5 = $d.Content
$findRange5.Find.Execute("6.3) Modelo da dispersão axial", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$programaPara2 = $findRange5.Paragraphs(1)
$programaPara2.Range.InsertParagraphAfter() | Out-Null
$newProgramaPara = $programaPara2.Next()
$newProgramaPara.Range.Text = "1. Introduction to Reactors: Basic concepts.2.Ideal Models of Isothermal Chemical Reactors - Simple Reactions: 2.1) Fundamental equations for reactor design. 2.2) Batch reactor (BSTR). 2.3) Continuous stirred-tank reactor (CSTR). 2.4) Plug-flow reactor (PFR). 2.5) Performance comparison of CSTR and PFR. 2.6) Cascade CSTR reactors. 2.7) Mixed association of reactors in series: CSTR and PFR. 2.8) Reactors with recycle. 2.9) Auto-catalytic reactions. 2.10) Semi-continuous reactors.3.Multiple Reactions in Ideal Reactors: 3.1) General concepts: optimization, yield, and selectivity. 3.2) Parallel reactions and series reactions.3.3) Systems with series-parallel reactions: multiple substitution reactions and polymerization reactions.3.4) Simple optimization problems.4.Thermal Effects in Ideal Reactors: 4.1) Energy balance equation. 4.2) Energy balance applied to BSTR. 4.3) Energy balance applied to CSTR. 4.4) Energy balance applied to PFR.5.Heterogeneous Catalytic Reactors: 5.1) Introduction. 5.2) Effect of physical processes on reaction rate:5.2.1 - Interfacial phenomena. 5.2.2 - Intraparticle phenomena. 5.2.3 - Diffusion and reaction in porous catalysts. 5.3) Calculation of fixed-bed reactors. 5.4) Three-phase reactors.6.Non-Ideal Reactors: 6.1) Residence time distribution. 6.2) Model for continuous stirred-tanks in series. 6.3) Axial dispersion model."
$italicRange3 = $programaPara2.Next().Range
$italicRange3.MoveEnd(1, -1)
$italicRange3.Font.Italic = $true

# ---------------------------------------------------------------------
# 8. Bibliografia: collapse into a single run, inserting the two new
#    section headers and the three new references.
# ---------------------------------------------------------------------
$findRange6 = $d.Content
$findRange6.Find.Execute("FOGLER, H. S. Elementos de Engenharia das Reações Químicas.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$biblioPara = $findRange6.Paragraphs(1)
$biblioPara.Range.Text = "Bibliografia Básica :FOGLER, H. S. Elementos de Engenharia das Reações Químicas. 3. ed. Rio de Janeiro: LTC Editora, 2002.LEVENSPIEL, O. Chemical Reaction Engineering. 3. ed. New York: John Wiley & Sons, 1998.HILL, C.G. An Introduction to Chemical Engineering Kinetics and Reactor Design. New York: John Wiley&Sons, 1977.Bibliografia Complementar:SMITH, J.M. Chemical Engineering Kinetics. 3rd. ed. New York :  McGraw-Hill, 1981.DENBIGH, K.; TURNER, R. Introduction to Chemical Reaction Design. Cambridge: Cambridge University Press, 1970.FROMENT, G.F.; BISCHOFF, K.B. Chemical Reactor Analysis And Design. 2nd ed.  New York: John Wiley & Sons, 1990.AGUIAR, L. G. Problemas de cinética e reatores químicos. Curitiba: Appris Editora, 2023.VAN SANTEN, R.A.; Niemantsverdriet, J.W. Chemical kinetics and catalysis. New York: Plenum Press, 1995.Missen, R.W.; Mims, C.A.; Saville, B.A. Introduction to chemical reaction engineering and kinetics. New York: J. Wiley, 1999.Rothenberg, G. Catalysis: concepts and green applications. Weinheim: Wiley-VCH, 2008 Chichester.Salmi, T.O.; Mikkola, J.; Warna, J.P. Chemical reaction engineering and reactor technology. Boca Raton: CRC Press/Taylor & Francis, 2011."

Write-Output "All edits applied"
